$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.17831799999999
$ws.Range("H2").Value = 288.534954
$ws.Range("I2").Value = 0.7237598617297997
$ws.Range("J2").Value = 0.7237598617297996
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 693.8410297858586
$ws.Range("R2").Value = 6244.569268072727
$ws.Range("S2").Value = 0.3393582462254739
$ws.Range("T2").Value = 0.3393582462254739

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.17831799999999
$ws.Range("H3").Value = 288.534954
$ws.Range("I3").Value = 0.7237598617297997
$ws.Range("J3").Value = 0.7237598617297996
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 683.8615675101786
$ws.Range("R3").Value = 6154.754107591607
$ws.Range("S3").Value = 0.3344772826174364
$ws.Range("T3").Value = 0.3344772826174364

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.17831799999999
$ws.Range("H4").Value = 288.534954
$ws.Range("I4").Value = 0.7237598617297997
$ws.Range("J4").Value = 0.7237598617297996
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 102.0736962395673
$ws.Range("R4").Value = 918.6632661561058
$ws.Range("S4").Value = 0.04992433288688936
$ws.Range("T4").Value = 0.04992433288688935

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 13.23504133333333
$ws.Range("H5").Value = 39.705124
$ws.Range("I5").Value = 0.09959616558694152
$ws.Range("J5").Value = 0.0995961655869415
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 95.47905285657421
$ws.Range("R5").Value = 859.3114757091679
$ws.Range("S5").Value = 0.04669888711925341
$ws.Range("T5").Value = 0.0466988871192534

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.23504133333333
$ws.Range("H6").Value = 39.705124
$ws.Range("I6").Value = 0.09959616558694152
$ws.Range("J6").Value = 0.0995961655869415
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 94.10578496782755
$ws.Range("R6").Value = 846.952064710448
$ws.Range("S6").Value = 0.0460272206101877
$ws.Range("T6").Value = 0.04602722061018769

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.23504133333333
$ws.Range("H7").Value = 39.705124
$ws.Range("I7").Value = 0.09959616558694152
$ws.Range("J7").Value = 0.0995961655869415
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 14.04630083858178
$ws.Range("R7").Value = 126.416707547236
$ws.Range("S7").Value = 0.006870057857500412
$ws.Range("T7").Value = 0.006870057857500411

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 23.47369766666667
$ws.Range("H8").Value = 70.421093
$ws.Range("I8").Value = 0.1766439726832589
$ws.Range("J8").Value = 0.1766439726832589
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 169.3418527232084
$ws.Range("R8").Value = 1524.076674508876
$ws.Range("S8").Value = 0.08282524625339155
$ws.Range("T8").Value = 0.08282524625339154

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 23.47369766666667
$ws.Range("H9").Value = 70.421093
$ws.Range("I9").Value = 0.1766439726832589
$ws.Range("J9").Value = 0.1766439726832589
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 166.9062218533151
$ws.Range("R9").Value = 1502.155996679836
$ws.Range("S9").Value = 0.08163397709377623
$ws.Range("T9").Value = 0.08163397709377622

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 23.47369766666667
$ws.Range("H10").Value = 70.421093
$ws.Range("I10").Value = 0.1766439726832589
$ws.Range("J10").Value = 0.1766439726832589
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 24.91254926340855
$ws.Range("R10").Value = 224.212943370677
$ws.Range("S10").Value = 0.01218474933609116
$ws.Range("T10").Value = 0.01218474933609116

